$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.418.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.842.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6270"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07405"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.839.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6701"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001038"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.256"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.383.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "234.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.470"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07274"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.485"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.481"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.034"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.040"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.161"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.817"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7145"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.574"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01839"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.781"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.232.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.796"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9554"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.01%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.997.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000119"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.16%  "
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.959"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.909"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3885"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.83%  "
